# Generate Report for Handoff
# Updates the "Priority" column (ht) for newly-flagged handoff rows and
# refreshes the related "Latest Handoff/HO Xliff Generate" timestamps
# across the Overview / zh-cn / de-de worksheets.

$wb = $excel.ActiveWorkbook

$ovWs   = $wb.Worksheets.Item("Overview")
$zhWs   = $wb.Worksheets.Item("zh-cn")
$deWs   = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Priority column (E) on zh-cn and de-de flips from blank to "ht"
    $zhWs.Range("E$r").Value = "ht"
    $deWs.Range("E$r").Value = "ht"

    # Latest Handoff Datetime (H) bumped by the new handoff generation
    $zhWs.Range("H$r").Value = "2016-09-04 10:24:45"
    $deWs.Range("H$r").Value = "2016-09-04 10:24:50"

    # Overview sheet's Latest HO Xliff Generate Date (G) mirrors de-de's bump
    $ovWs.Range("G$r").Value = "2016-09-04 10:24:50"
}
